$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update column F (dSF) values for rows 2-9 per repull/push of data
$ws.Range("F2").Value = -7
$ws.Range("F3").Value = -2
$ws.Range("F4").Value = 1
$ws.Range("F5").Value = -3
$ws.Range("F6").Value = -1
$ws.Range("F8").Value = -8
$ws.Range("F9").Value = 2
